$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn the "not started" row 25 (Make separate linked lists...) into a
#     "done" (green) row like row 26, picking up its number format / fill
#     and leaving the (now blank) D/E cells formatted the same way. ---
$ws.Range("A26:E26").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Replace the old "Loop Unrolling" placeholder row (27) with the new
#     "Dead code" entry, matching the yellow (in-progress) formatting used
#     by row 20. ---
$ws.Range("A20:E20").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A27").Value = "?"
$ws.Range("B27").Value = "Me "
$ws.Range("C27").Value = "Dead code"

# --- Fix a typo in the existing "Dead store" row (row 26) ---
$ws.Range("E26").Value = "Removes and warns"

$ws.Range("D27").Value = "Y"
$ws.Range("E27").Value = "Needs testing. Should work."

# --- Append a brand new row 28, styled like the other "done" rows but with
#     the fill cleared back out (no colour), mirroring rows 24/25. ---
$ws.Range("A26:E26").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A28:E28").Interior.ColorIndex = -4142   # xlNone

$ws.Range("A28").Value = "?"
$ws.Range("B28").Value = "Me"
$ws.Range("C28").Value = "Implement constant variables i.e. evaluate variable usage when determing constants (extra hard)"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""

# --- Update the view state to match the saved workbook (scrolled down a bit,
#     selection parked on the new last row). ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C31").Select()
